# daily auto push: 2026-01-08 06:50 UTC
# A new observation row is inserted at row 606 (date 2026/01/08, 木, 13, 165),
# pushing the existing rows 606-647 down to 607-648.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row above the current row 606, shifting rows 606:647 -> 607:648.
$ws.Rows.Item(606).Insert()

# Populate the newly inserted row. The leading apostrophe forces the
# date-like string to be stored as literal text (matching every other
# date cell in column A), rather than being auto-converted to a date
# serial number. ClearFormats() afterwards drops the "quote prefix"
# cell style that Excel applies when the apostrophe trick is used, so
# the cell ends up with the same (default) styling as its neighbours.
$ws.Range("A606").Value = "'2026/01/08"
$ws.Range("A606").ClearFormats()
$ws.Range("B606").Value = "木"
$ws.Range("C606").Value = 13
$ws.Range("D606").Value = 165
